# "Generate Report for Archive"
#
# Replaces every occurrence of the status text "Ready for handoff" with
# "In Translation" across all sheets, and shrinks the now-narrower
# "Status" columns to match (Overview!E:F, zh-cn!C, de-de!C).

$wb = $excel.ActiveWorkbook

$newText = "In Translation"

# --- Overview sheet: Status columns are E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = $newText
# Narrow the columns to reflect the shorter status text.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: Status column is C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C3").Value = $newText
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: Status column is C ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C3").Value = $newText
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
